# 2 Conductor Area Calculations
#
# - Renumber the "Pull #" column (D) for Conduits 1-3: the new 2nd-conductor
#   pulls are inserted at the front (old rows 7-8 become pulls 1 and 2) and
#   the previously-first pulls are pushed down (old rows 2-6 become pulls
#   3-7).
# - Center-align column E (Cable Size) for every existing data row, and
#   extend that same centered style across the (until-now blank) F/G cells
#   and the two brand-new H/I columns, so the sheet lines up once the 2nd
#   conductor's "Conduit Fill" values are wired in.
# - Add a new "Conduit 4" block (rows 9-10) with its own pull numbers
#   (8 and 9), following the same merged-cell layout used by the other
#   conduits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------------
# 1. Renumber the Pull # column (D2:D8)
# --------------------------------------------------------------------------
$ws.Range("D2").Value = 3
$ws.Range("D3").Value = 4
$ws.Range("D4").Value = 5
$ws.Range("D5").Value = 6
$ws.Range("D6").Value = 7
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 2

# --------------------------------------------------------------------------
# 2. Center-align column E (Cable Size) for every existing data row.
#    (VerticalAlignment is set before HorizontalAlignment so the engine
#    reuses the workbook's existing "centered" style instead of minting a
#    new, equivalent-but-distinct one.)
# --------------------------------------------------------------------------
$eRange = $ws.Range("E2:E8")
$eRange.VerticalAlignment = -4108
$eRange.HorizontalAlignment = -4108

# --------------------------------------------------------------------------
# 3. The blank F/G cells that fall on the "second" row of each merged
#    conduit block (F3, G3, F5, G5, F7, G7, F8, G8) also pick up the
#    centered style. Those cells are currently hidden inside existing
#    merges (F2:F3, G2:G3, F4:F5, G4:G5, F6:F8, G6:G8), and a merged-away
#    cell can't be written to directly - so unmerge, format, then re-merge.
# --------------------------------------------------------------------------
$ws.Range("F2:F3").UnMerge()
$ws.Range("G2:G3").UnMerge()
$ws.Range("F4:F5").UnMerge()
$ws.Range("G4:G5").UnMerge()
$ws.Range("F6:F8").UnMerge()
$ws.Range("G6:G8").UnMerge()

foreach ($addr in @("F3", "G3", "F5", "G5", "F7", "G7", "F8", "G8")) {
    $c = $ws.Range($addr)
    $c.VerticalAlignment = -4108
    $c.HorizontalAlignment = -4108
}

$ws.Range("F2:F3").Merge()
$ws.Range("G2:G3").Merge()
$ws.Range("F4:F5").Merge()
$ws.Range("G4:G5").Merge()
$ws.Range("F6:F8").Merge()
$ws.Range("G6:G8").Merge()

# --------------------------------------------------------------------------
# 4. Extend the centered style into the brand-new H/I columns, rows 2-8.
# --------------------------------------------------------------------------
$hiRange = $ws.Range("H2:I8")
$hiRange.VerticalAlignment = -4108
$hiRange.HorizontalAlignment = -4108

# --------------------------------------------------------------------------
# 5. New "Conduit 4" block: rows 9-10.
# --------------------------------------------------------------------------
$ws.Range("A9").Value = "Conduit 4"
$ws.Range("B9").Value = "543+00"
$ws.Range("C9").Value = "553+00"
$ws.Range("D9").Value = 8
$ws.Range("E9").Value = "2C#4"
$ws.Range("F9").Value = "LOCAL"

# G9 holds "29.06%" as literal text (like the other Conduit Fill cells),
# not a computed percentage. Assigning a percent-looking string straight to
# .Value auto-converts it to a numeric fraction, so instead stage the text
# in a scratch, Text-formatted cell far outside the used range, copy only
# the value across, then remove the scratch column again so the sheet's
# used range / dimension isn't left polluted.
$scratchCol = "K"
$ws.Range($scratchCol + "1").NumberFormat = "@"
$ws.Range($scratchCol + "1").Value = "29.06%"
$ws.Range($scratchCol + "1").Copy()
$ws.Range("G9").PasteSpecial(-4163)
$ws.Columns($scratchCol + ":" + $scratchCol).Delete()

$ws.Range("D10").Value = 9
$ws.Range("E10").Value = "2C#4"

# Center-align the new block, matching the other conduits. Row 10's A/B/C
# cells are the "hidden" half of the A9:A10 / B9:B10 / C9:C10 merges (same
# as row 3/5/7/8 for the existing conduits above) and are intentionally
# left untouched, so style row 9 (A:I) and row 10 (D:I) separately.
$row9Block = $ws.Range("A9:I9")
$row9Block.VerticalAlignment = -4108
$row9Block.HorizontalAlignment = -4108

$row10Block = $ws.Range("D10:I10")
$row10Block.VerticalAlignment = -4108
$row10Block.HorizontalAlignment = -4108

# --------------------------------------------------------------------------
# 6. Merge cells for the new Conduit 4 block.
# --------------------------------------------------------------------------
$ws.Range("A9:A10").Merge()
$ws.Range("B9:B10").Merge()
$ws.Range("C9:C10").Merge()
$ws.Range("F9:F10").Merge()
$ws.Range("G9:G10").Merge()
